$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Outbreak_Locations": add scenario 4 rows (copy of scenario 3
# pattern with scenario number bumped to 4)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Outbreak_Locations")

$ws1.Cells.Item(20, 1).Value = 4
$ws1.Cells.Item(20, 2).Value = 850
$ws1.Cells.Item(20, 3).Value = 250

$ws1.Cells.Item(21, 1).Value = 4
$ws1.Cells.Item(21, 2).Value = 650
$ws1.Cells.Item(21, 3).Value = 850

$ws1.Cells.Item(22, 1).Value = 4
$ws1.Cells.Item(22, 2).Value = 450
$ws1.Cells.Item(22, 3).Value = 150

$ws1.Cells.Item(23, 1).Value = 4
$ws1.Cells.Item(23, 2).Value = 750
$ws1.Cells.Item(23, 3).Value = 250

# ---------------------------------------------------------------
# Sheet "Store_Locations": add scenario 4 rows (copy of scenario 3
# store rows with scenario number bumped to 4)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Store_Locations")

$ws2.Cells.Item(11, 1).Value = 4
$ws2.Cells.Item(11, 2).Value = 112
$ws2.Cells.Item(11, 3).Value = 198
$ws2.Cells.Item(11, 4).Value = "Chain 1"

$ws2.Cells.Item(12, 1).Value = 4
$ws2.Cells.Item(12, 2).Value = 105
$ws2.Cells.Item(12, 3).Value = 855
$ws2.Cells.Item(12, 4).Value = "Chain 1"

# ---------------------------------------------------------------
# Sheet "Population": introduce a "population_type" column and a
# new scenario 4 entry ("random" type, no fixed population_per_cell)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Population")

# Write the (moved) population_per_cell column into C first so the
# existing text keeps referring to the same shared string.
$ws3.Cells.Item(1, 3).Value = "population_per_cell"
$ws3.Cells.Item(2, 3).Value = 5
$ws3.Cells.Item(3, 3).Value = 5
$ws3.Cells.Item(4, 3).Value = 5

# Now repurpose column B as the new population_type column.
$ws3.Cells.Item(1, 2).Value = "population_type"
$ws3.Cells.Item(2, 2).Value = "uniform"
$ws3.Cells.Item(3, 2).Value = "uniform"
$ws3.Cells.Item(4, 2).Value = "uniform"

# New scenario 4 row: random population, no population_per_cell value.
$ws3.Cells.Item(5, 1).Value = 4
$ws3.Cells.Item(5, 2).Value = "random"

# Match column C's width to the new header text.
$ws3.Columns.Item(3).ColumnWidth = 17.44140625

# ---------------------------------------------------------------
# View / selection state: Store_Locations and Population get a new
# selected cell, and Outbreak_Locations ends up the active sheet with
# H22 selected (scrolled down so row 7 is at the top).
# ---------------------------------------------------------------
$ws2.Range("A13").Select()
$ws3.Range("A5").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("H22").Select()

Write-Host "done"
